$d = $word.ActiveDocument

# Each (old -> new) text replacement from the diff, applied via Find/Replace
# on the whole document story. Each "old" string is unique in the document,
# so MatchWholeWord / non-wildcard exact matches are sufficient and safe.
$replacements = @(
    @("2026-01-15 Thursday", "2026-01-16 Friday"),
    @("742÷8=92, 6", "370÷8=46, 2"),
    @("241÷2=120, 1", "464÷7=66, 2"),
    @("590÷2=295, 0", "489÷3=163, 0"),
    @("343÷3=114, 1", "914÷6=152, 2"),
    @("606÷6=101, 0", "865÷7=123, 4"),
    @("182÷6=30, 2", "274÷4=68, 2"),
    @("628÷2=314, 0", "685÷4=171, 1"),
    @("609÷5=121, 4", "384÷6=64, 0"),
    @("888÷7=126, 6", "207÷2=103, 1"),
    @("975÷9=108, 3", "816÷7=116, 4"),
    @("454÷9=50, 4", "605÷5=121, 0"),
    @("515÷8=64, 3", "853÷6=142, 1"),
    @("417÷6=69, 3", "217÷9=24, 1"),
    @("397÷6=66, 1", "662÷8=82, 6"),
    @("992÷8=124, 0", "828÷7=118, 2"),
    @("102÷4=25, 2", "483÷4=120, 3"),
    @("143÷4=35, 3", "107÷3=35, 2"),
    @("809÷4=202, 1", "455÷5=91, 0"),
    @("798÷8=99, 6", "968÷2=484, 0"),
    @("312÷5=62, 2", "528÷8=66, 0"),
    @("864÷4=216, 0", "394÷8=49, 2"),
    @("135÷9=15, 0", "230÷7=32, 6"),
    @("508÷3=169, 1", "804÷2=402, 0"),
    @("649÷4=162, 1", "154÷7=22, 0"),
    @("312÷2=156, 0", "402÷8=50, 2"),
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
